# Update column F (dSF) values per repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -2
    4  = -4
    5  = -1
    6  = -1
    9  = -6
    12 = -3
    13 = -5
    14 = 24
    15 = 5
    16 = 4
    18 = 1
    19 = 2
    20 = -2
    21 = -1
    22 = 6
    23 = 2
    24 = -1
    26 = -2
    28 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
